# Kedar Jadhav.xlsx — add ownTeam/oppTeam columns and append new match rows.
#
# Before: A:venue B:date C:result D:batsman E:totalRuns F:totalBalls
#         G:total4s H:total6s I:sr   (1 data row)
# After:  A:venue B:date C:result D:ownTeam E:oppTeam F:batsman G:totalRuns
#         H:totalBalls I:total4s J:total6s K:sr   (5 data rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("venue","date","result","ownTeam","oppTeam","batsman","totalRuns","totalBalls","total4s","total6s","sr")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Columns G..K (totalRuns, totalBalls, total4s, total6s, sr) hold numeric-looking
# text; prefix with an apostrophe so Excel keeps them as text, matching the
# workbook's existing "number stored as text" data.
$data = @(
    @(" Abu Dhabi", " October 07 2020", "KKR won by 10 runs", "Chennai Super Kings", "Kolkata Knight Riders", "Kedar Jadhav ", "'7", "'12", "'1", "'0", "'58.33"),
    @(" Dubai (DSC)", " September 25 2020", "Capitals won by 44 runs", "Chennai Super Kings", "Delhi Capitals", "Kedar Jadhav ", "'26", "'21", "'3", "'0", "'123.80"),
    @(" Sharjah", " September 22 2020", "Royals won by 16 runs", "Chennai Super Kings", "Rajasthan Royals", "Kedar Jadhav ", "'22", "'16", "'3", "'0", "'137.50"),
    @(" Abu Dhabi", " October 19 2020", "Royals won by 7 wickets (with 15 balls remaining)", "Chennai Super Kings", "Rajasthan Royals", "Kedar Jadhav ", "'4", "'7", "'0", "'0", "'57.14"),
    @(" Dubai (DSC)", " October 02 2020", "Sunrisers won by 7 runs", "Chennai Super Kings", "Sunrisers Hyderabad", "Kedar Jadhav ", "'3", "'10", "'0", "'0", "'30.00")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
